# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# for the cryptos worksheet per the commit's refreshed data snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.912.84'
$ws.Range('E2').Value = '  +1.39%  '
$ws.Range('D3').Value = '3.388.56'
$ws.Range('E3').Value = '  +1.26%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '559.49'
$ws.Range('E5').Value = '  +1.49%  '
$ws.Range('D6').Value = '174.91'
$ws.Range('E6').Value = '  +0.98%  '
$ws.Range('E7').Value = '  +1.85%  '
$ws.Range('D8').Value = '3.378.73'
$ws.Range('E8').Value = '  +1.28%  '
$ws.Range('E9').Value = '  -0.06%  '
$ws.Range('E10').Value = '  +11.38%  '
$ws.Range('E11').Value = '  +2.71%  '
$ws.Range('D12').Value = '54.78'
$ws.Range('E12').Value = '  +1.62%  '
$ws.Range('E13').Value = '  +5.13%  '
$ws.Range('E14').Value = '  +2.55%  '
$ws.Range('D15').Value = '3.930.01'
$ws.Range('E15').Value = '  +4.05%  '
$ws.Range('D16').Value = '18.34'
$ws.Range('E16').Value = '  +0.15%  '
$ws.Range('E17').Value = '  +1.78%  '
$ws.Range('D18').Value = '3.384.99'
$ws.Range('E18').Value = '  +1.53%  '
$ws.Range('D19').Value = '64.903.96'
$ws.Range('E19').Value = '  +1.38%  '
$ws.Range('D20').Value = '11.83'
$ws.Range('E20').Value = '  +0.74%  '
$ws.Range('D21').Value = '0.994'
$ws.Range('E21').Value = '  +1.53%  '
$ws.Range('D22').Value = '474.43'
$ws.Range('E22').Value = '  +15.48%  '
$ws.Range('D23').Value = '4.98'
$ws.Range('E23').Value = '  +13.40%  '
$ws.Range('E24').Value = '  +1.84%  '
$ws.Range('D25').Value = '87.06'
$ws.Range('E25').Value = '  +4.69%  '
$ws.Range('D26').Value = '13.51'
$ws.Range('E26').Value = '  -2.56%  '
$ws.Range('E27').Value = '  +6.54%  '
$ws.Range('D28').Value = '10.88'
$ws.Range('E28').Value = '  +2.73%  '
$ws.Range('D29').Value = '8.77'
$ws.Range('E29').Value = '  +1.28%  '
$ws.Range('D30').Value = '31.07'
$ws.Range('E30').Value = '  +6.60%  '
$ws.Range('D31').Value = '6.68'
$ws.Range('E31').Value = '  +3.59%  '
$ws.Range('D32').Value = '11.54'
$ws.Range('E32').Value = '  +1.27%  '
$ws.Range('D33').Value = '61.72'
$ws.Range('E33').Value = '  +6.45%  '
$ws.Range('D34').Value = '571.47'
$ws.Range('E34').Value = '  -1.63%  '
$ws.Range('E35').Value = '  +1.43%  '
$ws.Range('E36').Value = '  +0.06%  '
$ws.Range('D37').Value = '3.56'
$ws.Range('E37').Value = '  +3.76%  '
$ws.Range('E38').Value = '  -5.40%  '
$ws.Range('D39').Value = '35.72'
$ws.Range('E40').Value = '  +2.03%  '
$ws.Range('D41').Value = '0.371'
$ws.Range('E41').Value = '  +0.99%  '
$ws.Range('D42').Value = '3.091.01'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  +0.06%  '
$ws.Range('E44').Value = '  +1.92%  '
$ws.Range('D45').Value = '0.0416'
$ws.Range('E45').Value = '  +3.52%  '
$ws.Range('E46').Value = '  +5.69%  '
$ws.Range('E47').Value = '  +1.61%  '
$ws.Range('D48').Value = '3.15'
$ws.Range('E48').Value = '  -2.61%  '
$ws.Range('E49').Value = '  -0.17%  '
$ws.Range('D50').Value = '137.49'
$ws.Range('E50').Value = '  +2.12%  '
$ws.Range('E51').Value = '  +3.43%  '
